# Atualizacao de bases das ligas (Lithuania A Lyga): corrige a ordem de
# duas partidas existentes (id 99 / id 101), atualiza os dados de uma
# partida que ainda nao tinha resultado (id 106) e adiciona uma nova
# partida (id 107, Hegelmann Litauen x FK Transinvest).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: Swap data between row 101 and row 103 (columns B:AC; column A id stays put) ---
# Row 101 (after) values:
$ws.Range("B101").Value = 6732836
$ws.Range("C101").Value = "Lithuania A Lyga"
$ws.Range("D101").Value = "Lithuania A Lyga"
$ws.Range("E101").Value = 45242.41319444445
$ws.Range("F101").Value = "FK Siauliai"
$ws.Range("G101").Value = "Banga Gargzdai"
$ws.Range("H101").Value = 3
$ws.Range("I101").Value = 0
$ws.Range("J101").Value = "H"
$ws.Range("K101").Value = 1.222
$ws.Range("L101").Value = 5.5
$ws.Range("M101").Value = 9
$ws.Range("N101").Value = 1.363
$ws.Range("O101").Value = 4.5
$ws.Range("P101").Value = 7
$ws.Range("Q101").Value = -1.25
$ws.Range("R101").Value = 1.9
$ws.Range("S101").Value = 1.9
$ws.Range("T101").Value = 2.5
$ws.Range("U101").Value = 1.975
$ws.Range("V101").Value = 1.825
$ws.Range("W101").Value = 0.363
$ws.Range("X101").Value = -1
$ws.Range("Y101").Value = -1
$ws.Range("Z101").Value = 0.8999999999999999
$ws.Range("AA101").Value = -1
$ws.Range("AB101").Value = 0.9750000000000001
$ws.Range("AC101").Value = -1

# Row 103 (after) values:
$ws.Range("B103").Value = 7465686
$ws.Range("C103").Value = "Lithuania A Lyga"
$ws.Range("D103").Value = "Lithuania A Lyga"
$ws.Range("E103").Value = 45242.41319444445
$ws.Range("F103").Value = "FK Kauno Zalgiris"
$ws.Range("G103").Value = "Hegelmann Litauen"
$ws.Range("H103").Value = 4
$ws.Range("I103").Value = 2
$ws.Range("J103").Value = "H"
$ws.Range("K103").Value = 2.3
$ws.Range("L103").Value = 4
$ws.Range("M103").Value = 2.3
$ws.Range("N103").Value = 2.55
$ws.Range("O103").Value = 4
$ws.Range("P103").Value = 2.2
$ws.Range("Q103").Value = 0.25
$ws.Range("R103").Value = 1.8
$ws.Range("S103").Value = 2
$ws.Range("T103").Value = 2.75
$ws.Range("U103").Value = 1.85
$ws.Range("V103").Value = 1.95
$ws.Range("W103").Value = 1.55
$ws.Range("X103").Value = -1
$ws.Range("Y103").Value = -1
$ws.Range("Z103").Value = 0.8
$ws.Range("AA103").Value = -1
$ws.Range("AB103").Value = 0.8500000000000001
$ws.Range("AC103").Value = -1

# --- Step 2: Update row 108 (add H/I/J, update K:AC) ---
$ws.Range("B108").Value = 7862034
$ws.Range("C108").Value = "Lithuania A Lyga"
$ws.Range("D108").Value = "Lithuania A Lyga"
$ws.Range("E108").Value = 45354.33333333334
$ws.Range("F108").Value = "FK Dainava Alytus"
$ws.Range("G108").Value = "FK Zalgiris Vilnius"
$ws.Range("H108").Value = 0
$ws.Range("I108").Value = 1
$ws.Range("J108").Value = "A"
$ws.Range("K108").Value = 6.5
$ws.Range("L108").Value = 4
$ws.Range("M108").Value = 1.4
$ws.Range("N108").Value = 6.5
$ws.Range("O108").Value = 3.6
$ws.Range("P108").Value = 1.45
$ws.Range("Q108").Value = 1
$ws.Range("R108").Value = 1.95
$ws.Range("S108").Value = 1.85
$ws.Range("T108").Value = 2.25
$ws.Range("U108").Value = 2
$ws.Range("V108").Value = 1.8
$ws.Range("W108").Value = -1
$ws.Range("X108").Value = -1
$ws.Range("Y108").Value = 0.45
$ws.Range("Z108").Value = 0
$ws.Range("AA108").Value = 0
$ws.Range("AB108").Value = -1
$ws.Range("AC108").Value = 0.8

# --- Step 3: Add new row 109 ---
$ws.Range("A109").Value = 107
$ws.Range("B109").Value = 7862904
$ws.Range("C109").Value = "Lithuania A Lyga"
$ws.Range("D109").Value = "Lithuania A Lyga"
$ws.Range("E109").Value = 45354.5625
$ws.Range("F109").Value = "Hegelmann Litauen"
$ws.Range("G109").Value = "FK Transinvest"
$ws.Range("H109").Value = 2
$ws.Range("I109").Value = 3
$ws.Range("J109").Value = "A"
$ws.Range("K109").Value = 1.8
$ws.Range("L109").Value = 3.5
$ws.Range("M109").Value = 3.6
$ws.Range("N109").Value = 1.95
$ws.Range("O109").Value = 3.5
$ws.Range("P109").Value = 3.1
$ws.Range("Q109").Value = -0.25
$ws.Range("R109").Value = 1.775
$ws.Range("S109").Value = 2.025
$ws.Range("T109").Value = 2.5
$ws.Range("U109").Value = 1.9
$ws.Range("V109").Value = 1.9
$ws.Range("W109").Value = -1
$ws.Range("X109").Value = -1
$ws.Range("Y109").Value = 2.1
$ws.Range("Z109").Value = -1
$ws.Range("AA109").Value = 1.025
$ws.Range("AB109").Value = 0.8999999999999999
$ws.Range("AC109").Value = -1

# Apply style formatting to new row 109 (match row 108 formatting for A and E columns)
$ws.Range("A108").Copy()
$ws.Range("A109").PasteSpecial(-4122, $null, $false, $false)
$ws.Range("E108").Copy()
$ws.Range("E109").PasteSpecial(-4122, $null, $false, $false)
